$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 359 (shifts existing rows 359:432 down to 360:433)
$ws.Rows("359").Insert()

# Populate the newly inserted row with a new weekly data record
$ws.Range("A359").Value = 7
$ws.Range("B359").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C359").Value = "Ñuble"
$ws.Range("D359").Value = 45015
$ws.Range("E359").Value = 16
$ws.Range("F359").Value = 100112023
$ws.Range("G359").Value = "Brócoli"
$ws.Range("H359").Value = "Sin especificar"
$ws.Range("I359").Value = "Segunda"
$ws.Range("J359").Value = 100
$ws.Range("K359").Value = 1000
$ws.Range("L359").Value = 1000
$ws.Range("M359").Value = 1000
$ws.Range("N359").Value = "$/unidad"
$ws.Range("O359").Value = "Provincia de Diguillín"
$ws.Range("P359").Value = 1000
$ws.Range("Q359").Value = 1
$ws.Range("R359").Value = "Hortaliza"
